$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-7 (columns A through T), per the NATMI LR-pair output
# for Fgf17 -> Fgfr3 across FAPs / sCs / ECs clusters.
$data = @(
    @("FAPs", "Fgf17", "Fgfr3", "ECs", 1.0, 0.3333333333333333, 0.05925133333333333, 0.177754, 0.129311749608618, 0.129311749608618, 2.0, 0.6666666666666666, 4.220261333333333, 12.660784, 0.6739259863235564, 0.6739259863235564, 0.2500561110151111, 2.250504999136, 0.08714654839821262, 0.08714654839821262),
    @("FAPs", "Fgf17", "Fgfr3", "FAPs", 1.0, 0.3333333333333333, 0.05925133333333333, 0.177754, 0.129311749608618, 0.129311749608618, 3.0, 1.0, 1.296447666666667, 3.889343, 0.2070274097896007, 0.2070274097896007, 0.07681625284688888, 0.691346275622, 0.02677107657683359, 0.02677107657683359),
    @("FAPs", "Fgf17", "Fgfr3", "sCs", 1.0, 0.3333333333333333, 0.05925133333333333, 0.177754, 0.129311749608618, 0.129311749608618, 3.0, 1.0, 0.745494, 2.236482, 0.119046603886843, 0.119046603886843, 0.044171513492, 0.397543621428, 0.01539412463357177, 0.01539412463357177),
    @("sCs", "Fgf17", "Fgfr3", "ECs", 3.0, 1.0, 0.398954, 1.196862, 0.8706882503913821, 0.8706882503913821, 2.0, 0.6666666666666666, 4.220261333333333, 12.660784, 0.6739259863235564, 0.6739259863235564, 1.683690139978667, 15.153211259808, 0.5867794379253438, 0.5867794379253438),
    @("sCs", "Fgf17", "Fgfr3", "FAPs", 3.0, 1.0, 0.398954, 1.196862, 0.8706882503913821, 0.8706882503913821, 3.0, 1.0, 1.296447666666667, 3.889343, 0.2070274097896007, 0.2070274097896007, 0.5172229824073333, 4.655006841666, 0.1802563332127671, 0.1802563332127671),
    @("sCs", "Fgf17", "Fgfr3", "sCs", 3.0, 1.0, 0.398954, 1.196862, 0.8706882503913821, 0.8706882503913821, 3.0, 1.0, 0.745494, 2.236482, 0.119046603886843, 0.119046603886843, 0.297417813276, 2.676760319484, 0.1036524792532713, 0.1036524792532713)
)

$startRow = 2
$r = $startRow
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
